# Inserts two new data rows (new rows 507 and 508) into the "Crespo record"
# price table on Sheet1, pushing the previously-existing rows 507:613 down
# to 509:615 (dimension grows from A1:R613 to A1:R615).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 507:613 down by two rows, copying formatting from the
# row above (matches the style used for row 506 -> the row just inserted).
$ws.Rows("507:508").Insert()

# --- New row 507 -------------------------------------------------------
$ws.Range("A507").Value = 3
$ws.Range("B507").Value = "Femacal de La Calera"
$ws.Range("C507").Value = "Coquimbo"
$ws.Range("D507").Value = 44711
$ws.Range("E507").Value = 5
$ws.Range("F507").Value = 100112006
$ws.Range("G507").Value = "Repollo"
$ws.Range("H507").Value = "Crespo record"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 1900
$ws.Range("K507").Value = 1000
$ws.Range("L507").Value = 1100
$ws.Range("M507").Value = 1050
$ws.Range("N507").Value = "$/unidad"
$ws.Range("O507").Value = "Provincia de Quillota"
$ws.Range("P507").Value = 1050
$ws.Range("Q507").Value = 1
$ws.Range("R507").Value = "Hortaliza"

# --- New row 508 -------------------------------------------------------
$ws.Range("A508").Value = 3
$ws.Range("B508").Value = "Femacal de La Calera"
$ws.Range("C508").Value = "Coquimbo"
$ws.Range("D508").Value = 44711
$ws.Range("E508").Value = 5
$ws.Range("F508").Value = 100112006
$ws.Range("G508").Value = "Repollo"
$ws.Range("H508").Value = "Crespo record"
$ws.Range("I508").Value = "Segunda"
$ws.Range("J508").Value = 900
$ws.Range("K508").Value = 800
$ws.Range("L508").Value = 800
$ws.Range("M508").Value = 800
$ws.Range("N508").Value = "$/unidad"
$ws.Range("O508").Value = "Provincia de Quillota"
$ws.Range("P508").Value = 800
$ws.Range("Q508").Value = 1
$ws.Range("R508").Value = "Hortaliza"
